# Currency changes between dates sheet implementation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header typo "Simbulo" -> "Simbolo"
$ws.Range("D1").Value = "Símbolo"

# Update currency quote values
$ws.Range("B2").Value = 5.6
$ws.Range("B4").Value = 4.5
$ws.Range("B5").Value = 6.47
$ws.Range("B6").Value = 6.14
$ws.Range("B7").Value = 0.0492
$ws.Range("B8").Value = 7.57
$ws.Range("B10").Value = 0.0069
$ws.Range("B13").Value = 0.88

# Update report date in C16 while keeping it stored as text (its cell is
# formatted as a date, so a direct assignment would be auto-converted to a
# date serial number). Stage the text in a scratch cell formatted as Text,
# then paste only the value across, and remove the scratch column again.
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "04/11/2021"
$ws.Range("Z1").Copy()
$ws.Range("C16").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("Z1").EntireColumn.Delete()

# Update report time footer
$ws.Range("D16").Value = "21:57"

# Update the active selection to D5
$ws.Range("D5").Select()
